$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-4 first (P, I, D) so the shared-string table is built in the same
# order as the target file, then the header row (Max Temperature) last.
$ws.Range("A2").Value = "P"
$ws.Range("A3").Value = "I"
$ws.Range("A4").Value = "D"
$ws.Range("A1").Value = "Max Temperature"

$ws.Range("B1").Value = 50
$ws.Range("C1").Value = 100
$ws.Range("D1").Value = 150
$ws.Range("E1").Value = 200

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 8

$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 12

# Center-align the whole populated table (adds the 2nd cellXf used by every cell).
$ws.Range("A1:E4").HorizontalAlignment = -4108

# Widen column A to fit the longest label ("Max Temperature").
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the same cell selected as in the authored file.
$ws.Range("F9").Select() | Out-Null
